# Populate the daily shop-performance rows for 2026-02-xx (rows 93-101)
# with the newly uploaded figures. Column D ("总访客"? actually diff col D)
# already carries a shared formula (=Bn-Cn) and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    93  = @(6926, 5947, 112, 9, 9, 206, 5, 26)
    94  = @(14369, 13613, 128, 3, 3, 256, 3, 0)
    95  = @(4708, 4081, 121, 8, 8, 214, 7, 52)
    96  = @(4728, 4330, 99, 5, 5, 187, 3, 20)
    97  = @(414, 58, 9, 1, 4, 34, 0, 20)
    98  = @(4187, 3863, 113, 4, 5, 190, 5, 75)
    99  = @(5331, 4537, 102, 7, 5, 176, 8, 62)
    100 = @(5962, 5622, 86, 5, 5, 195, 10, 979)
    101 = @(5273, 4750, 92, 6, 6, 187, 4, 202)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    # D$row is left untouched - it already holds the shared formula =B-C
    # and will recalculate automatically.
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
    $ws.Range("H$row").Value = $vals[5]
    $ws.Range("I$row").Value = $vals[6]
    $ws.Range("J$row").Value = $vals[7]
}

# Restore the on-screen selection to match where the author ended up
# (frozen-pane view, bottom-right pane active cell K103).
$ws.Range("K103").Select()
